$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 433.33334
$ws.Cells.Item(18, 9).Value = 525
$ws.Cells.Item(18, 10).Value = 250
$ws.Cells.Item(18, 11).Value = 525
$ws.Cells.Item(18, 12).Value = 250
$ws.Cells.Item(18, 13).Value = -241
$ws.Cells.Item(18, 14).Value = -818

$ws.Cells.Item(43, 8).Value = 2499.5
$ws.Cells.Item(43, 9).Value = 2499.5
$ws.Cells.Item(43, 11).Value = 2499.5
$ws.Cells.Item(43, 13).Value = -2430.5

$ws.Cells.Item(112, 8).Value = 993
$ws.Cells.Item(112, 9).Value = 1110.6666
$ws.Cells.Item(112, 10).Value = 640
$ws.Cells.Item(112, 11).Value = 3331.9998
$ws.Cells.Item(112, 12).Value = 1920
$ws.Cells.Item(112, 13).Value = -2223.9998
$ws.Cells.Item(112, 14).Value = -4136

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4880.65
$ws.Cells.Item(2, 9).Value = 2922.8823
$ws.Cells.Item(2, 11).Value = 2922.8823
$ws.Cells.Item(2, 13).Value = -2809.8823

$ws.Cells.Item(4, 8).Value = 189.75
$ws.Cells.Item(4, 9).Value = 155
$ws.Cells.Item(4, 10).Value = 224.5
$ws.Cells.Item(4, 11).Value = 155
$ws.Cells.Item(4, 12).Value = 224.5
$ws.Cells.Item(4, 13).Value = -39
$ws.Cells.Item(4, 14).Value = -456.5

$ws.Cells.Item(32, 8).Value = 2570671
$ws.Cells.Item(32, 9).Value = 3352.7666
$ws.Cells.Item(32, 10).Value = 11128398
$ws.Cells.Item(32, 11).Value = 3352.7666
$ws.Cells.Item(32, 12).Value = 11128398
$ws.Cells.Item(32, 13).Value = -3065.7666
$ws.Cells.Item(32, 14).Value = -11128972

$ws.Cells.Item(45, 8).Value = 3958.889
$ws.Cells.Item(45, 9).Value = 2776
$ws.Cells.Item(45, 10).Value = 5437.5
$ws.Cells.Item(45, 11).Value = 2776
$ws.Cells.Item(45, 12).Value = 5437.5
$ws.Cells.Item(45, 13).Value = -2399
$ws.Cells.Item(45, 14).Value = -6191.5

$ws.Cells.Item(88, 8).Value = 1380.1111
$ws.Cells.Item(88, 9).Value = 395.25
$ws.Cells.Item(88, 10).Value = 2168
$ws.Cells.Item(88, 11).Value = 395.25
$ws.Cells.Item(88, 12).Value = 2168
$ws.Cells.Item(88, 13).Value = 10.75
$ws.Cells.Item(88, 14).Value = -2980

$ws.Cells.Item(91, 8).Value = 1380.1111
$ws.Cells.Item(91, 9).Value = 395.25
$ws.Cells.Item(91, 10).Value = 2168
$ws.Cells.Item(91, 11).Value = 395.25
$ws.Cells.Item(91, 12).Value = 2168
$ws.Cells.Item(91, 13).Value = 1008.75
$ws.Cells.Item(91, 14).Value = -4976

$ws.Cells.Item(110, 8).Value = 2797.5
$ws.Cells.Item(110, 9).Value = 2699.3333
$ws.Cells.Item(110, 11).Value = 2699.3333
$ws.Cells.Item(110, 13).Value = -654.3332999999998

$ws.Cells.Item(116, 8).Value = 4880.65
$ws.Cells.Item(116, 9).Value = 2922.8823
$ws.Cells.Item(116, 11).Value = 2922.8823
$ws.Cells.Item(116, 13).Value = -628.8823000000002

$ws.Cells.Item(122, 8).Value = 3746.625
$ws.Cells.Item(122, 9).Value = 3425.8572
$ws.Cells.Item(122, 11).Value = 10277.5716
$ws.Cells.Item(122, 13).Value = -7827.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4880.65
$ws.Cells.Item(3, 9).Value = 2922.8823
$ws.Cells.Item(3, 11).Value = 2922.8823
$ws.Cells.Item(3, 13).Value = -2808.8823

$ws.Cells.Item(76, 8).Value = 13099.667
$ws.Cells.Item(76, 10).Value = 13099.667
$ws.Cells.Item(76, 12).Value = 13099.667
$ws.Cells.Item(76, 14).Value = -13729.667

$ws.Cells.Item(79, 8).Value = 13099.667
$ws.Cells.Item(79, 10).Value = 13099.667
$ws.Cells.Item(79, 12).Value = 13099.667
$ws.Cells.Item(79, 14).Value = -15283.667

$ws.Cells.Item(88, 8).Value = 22000
$ws.Cells.Item(88, 10).Value = 22000
$ws.Cells.Item(88, 12).Value = 22000
$ws.Cells.Item(88, 14).Value = -22812

$ws.Cells.Item(91, 8).Value = 22000
$ws.Cells.Item(91, 10).Value = 22000
$ws.Cells.Item(91, 12).Value = 22000
$ws.Cells.Item(91, 14).Value = -24808

$ws.Cells.Item(134, 8).Value = 2794
$ws.Cells.Item(134, 9).Value = 2481.9443
$ws.Cells.Item(134, 10).Value = 4666.3335
$ws.Cells.Item(134, 11).Value = 7445.8329
$ws.Cells.Item(134, 12).Value = 13999.0005
$ws.Cells.Item(134, 13).Value = -4910.8329
$ws.Cells.Item(134, 14).Value = -19069.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1116.6666
$ws.Cells.Item(16, 10).Value = 1100
$ws.Cells.Item(16, 12).Value = 1100
$ws.Cells.Item(16, 14).Value = -1674

$ws.Cells.Item(50, 8).Value = 25000
$ws.Cells.Item(50, 10).Value = 30000
$ws.Cells.Item(50, 12).Value = 30000
$ws.Cells.Item(50, 14).Value = -31250

$ws.Cells.Item(51, 8).Value = 71296.664
$ws.Cells.Item(51, 10).Value = 71296.664
$ws.Cells.Item(51, 12).Value = 71296.664
$ws.Cells.Item(51, 14).Value = -72768.664

$ws.Cells.Item(58, 8).Value = 3202.0527
$ws.Cells.Item(58, 9).Value = 837.1667
$ws.Cells.Item(58, 10).Value = 4293.5386
$ws.Cells.Item(58, 11).Value = 837.1667
$ws.Cells.Item(58, 12).Value = 4293.5386
$ws.Cells.Item(58, 13).Value = -634.1667
$ws.Cells.Item(58, 14).Value = -4699.5386

$ws.Cells.Item(61, 8).Value = 71296.664
$ws.Cells.Item(61, 10).Value = 71296.664
$ws.Cells.Item(61, 12).Value = 71296.664
$ws.Cells.Item(61, 14).Value = -71992.664

$ws.Cells.Item(105, 8).Value = 3081
$ws.Cells.Item(105, 9).Value = 3081
$ws.Cells.Item(105, 11).Value = 3081
$ws.Cells.Item(105, 13).Value = -1334

$ws.Cells.Item(113, 8).Value = 1116.6666
$ws.Cells.Item(113, 10).Value = 1100
$ws.Cells.Item(113, 12).Value = 1100
$ws.Cells.Item(113, 14).Value = -5440

$ws.Cells.Item(136, 8).Value = 3202.0527
$ws.Cells.Item(136, 9).Value = 837.1667
$ws.Cells.Item(136, 10).Value = 4293.5386
$ws.Cells.Item(136, 11).Value = 2511.5001
$ws.Cells.Item(136, 12).Value = 12880.6158
$ws.Cells.Item(136, 13).Value = 38.4998999999998
$ws.Cells.Item(136, 14).Value = -17980.6158

$ws.Cells.Item(141, 8).Value = 84826
$ws.Cells.Item(141, 10).Value = 84826
$ws.Cells.Item(141, 12).Value = 84826
$ws.Cells.Item(141, 14).Value = -95186

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 33402082
$ws.Cells.Item(4, 9).Value = 55667068
$ws.Cells.Item(4, 10).Value = 4601.5
$ws.Cells.Item(4, 11).Value = 167001204
$ws.Cells.Item(4, 12).Value = 13804.5
$ws.Cells.Item(4, 13).Value = -167001092
$ws.Cells.Item(4, 14).Value = -14028.5

$ws.Cells.Item(5, 8).Value = 1724.75
$ws.Cells.Item(5, 10).Value = 1694.75
$ws.Cells.Item(5, 12).Value = 5084.25
$ws.Cells.Item(5, 14).Value = -5308.25

$ws.Cells.Item(135, 8).Value = 1724.75
$ws.Cells.Item(135, 10).Value = 1694.75
$ws.Cells.Item(135, 12).Value = 15252.75
$ws.Cells.Item(135, 14).Value = -20322.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 62561.312
$ws.Cells.Item(2, 9).Value = 142898.14
$ws.Cells.Item(2, 10).Value = 77.111115
$ws.Cells.Item(2, 11).Value = 142898.14
$ws.Cells.Item(2, 12).Value = 77.111115
$ws.Cells.Item(2, 13).Value = -142785.14
$ws.Cells.Item(2, 14).Value = -303.111115

$ws.Cells.Item(113, 8).Value = 8071.9287
$ws.Cells.Item(113, 9).Value = 5501.6665
$ws.Cells.Item(113, 11).Value = 5501.6665
$ws.Cells.Item(113, 13).Value = -3331.6665

$ws.Cells.Item(132, 8).Value = 36542.875
$ws.Cells.Item(132, 9).Value = 48481.74
$ws.Cells.Item(132, 10).Value = 6032.4443
$ws.Cells.Item(132, 11).Value = 145445.22
$ws.Cells.Item(132, 12).Value = 18097.3329
$ws.Cells.Item(132, 13).Value = -142915.22
$ws.Cells.Item(132, 14).Value = -23157.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7267.5
$ws.Cells.Item(7, 9).Value = 5173.8335
$ws.Cells.Item(7, 11).Value = 5173.8335
$ws.Cells.Item(7, 13).Value = -5061.8335

$ws.Cells.Item(40, 8).Value = 5238.231
$ws.Cells.Item(40, 9).Value = 2607.375
$ws.Cells.Item(40, 10).Value = 9447.6
$ws.Cells.Item(40, 11).Value = 2607.375
$ws.Cells.Item(40, 12).Value = 9447.6
$ws.Cells.Item(40, 13).Value = -2471.375
$ws.Cells.Item(40, 14).Value = -9719.6

$ws.Cells.Item(46, 8).Value = 4937.852
$ws.Cells.Item(46, 10).Value = 5236.579
$ws.Cells.Item(46, 12).Value = 5236.579
$ws.Cells.Item(46, 14).Value = -5612.579

$ws.Cells.Item(126, 8).Value = 7267.5
$ws.Cells.Item(126, 9).Value = 5173.8335
$ws.Cells.Item(126, 11).Value = 15521.5005
$ws.Cells.Item(126, 13).Value = -13051.5005

$ws.Cells.Item(132, 8).Value = 3883
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 14).Value = -14060

$ws.Cells.Item(136, 8).Value = 3422.3333
$ws.Cells.Item(136, 9).Value = 3422.3333
$ws.Cells.Item(136, 11).Value = 10266.9999
$ws.Cells.Item(136, 13).Value = -7716.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 12222856
$ws.Cells.Item(5, 9).Value = 15000033
$ws.Cells.Item(5, 10).Value = 6668500
$ws.Cells.Item(5, 11).Value = 15000033
$ws.Cells.Item(5, 12).Value = 6668500
$ws.Cells.Item(5, 13).Value = -14999921
$ws.Cells.Item(5, 14).Value = -6668724

$ws.Cells.Item(14, 8).Value = 763
$ws.Cells.Item(14, 9).Value = 796.875
$ws.Cells.Item(14, 10).Value = 724.2857
$ws.Cells.Item(14, 11).Value = 796.875
$ws.Cells.Item(14, 12).Value = 724.2857
$ws.Cells.Item(14, 13).Value = -628.875
$ws.Cells.Item(14, 14).Value = -1060.2857

$ws.Cells.Item(23, 8).Value = 2097.5
$ws.Cells.Item(23, 9).Value = 2097.5
$ws.Cells.Item(23, 11).Value = 2097.5
$ws.Cells.Item(23, 13).Value = -1868.5

$ws.Cells.Item(39, 8).Value = 99999
$ws.Cells.Item(39, 9).Value = 99999
$ws.Cells.Item(39, 11).Value = 99999
$ws.Cells.Item(39, 13).Value = -99586

$ws.Cells.Item(81, 8).Value = 857.4286
$ws.Cells.Item(81, 9).Value = 833.3333
$ws.Cells.Item(81, 10).Value = 1002
$ws.Cells.Item(81, 11).Value = 1666.6666
$ws.Cells.Item(81, 12).Value = 2004
$ws.Cells.Item(81, 13).Value = -605.6666
$ws.Cells.Item(81, 14).Value = -4126

$ws.Cells.Item(84, 8).Value = 857.4286
$ws.Cells.Item(84, 9).Value = 833.3333
$ws.Cells.Item(84, 10).Value = 1002
$ws.Cells.Item(84, 11).Value = 8333.333000000001
$ws.Cells.Item(84, 12).Value = 10020
$ws.Cells.Item(84, 13).Value = -3029.333000000001
$ws.Cells.Item(84, 14).Value = -20628

$ws.Cells.Item(122, 8).Value = 2101.75
$ws.Cells.Item(122, 9).Value = 2065.5454
$ws.Cells.Item(122, 11).Value = 6196.6362
$ws.Cells.Item(122, 13).Value = -3746.6362

$ws.Cells.Item(124, 8).Value = 20429
$ws.Cells.Item(124, 10).Value = 20429
$ws.Cells.Item(124, 12).Value = 20429
$ws.Cells.Item(124, 14).Value = -30249

$ws.Cells.Item(132, 8).Value = 3445.3
$ws.Cells.Item(132, 9).Value = 2030
$ws.Cells.Item(132, 10).Value = 4860.6
$ws.Cells.Item(132, 11).Value = 6090
$ws.Cells.Item(132, 12).Value = 14581.8
$ws.Cells.Item(132, 13).Value = -3560
$ws.Cells.Item(132, 14).Value = -19641.8
